$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Excel auto-converts plain "NN%" strings into percentage numbers when
# assigned via Range.Value, which would also silently swap the General
# number format for a Percent one (changing the cell style index).
# To keep these values as literal text (matching the source data, which
# stores them as plain strings like "55%") and keep the original style
# untouched, build the text in an off-sheet scratch cell via a formula
# (which always yields a text-typed result) and paste only the value.
$xlPasteValues = -4163

$ws.Range("E2").Value = "2026-02-26 04:18:38"
$ws.Range("O2").Value = "1.8 °C"
$ws.Range("E3").Value = "2026-02-26 04:18:41"
$ws.Range("Z1").Formula = "=""55%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H3").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("O3").Value = "1.4 °C"
$ws.Range("E4").Value = "2026-02-26 04:18:43"
$ws.Range("Z1").Formula = "=""93%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H4").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("K4").Value = "-0.1 MJ/m2"
$ws.Range("L4").Value = "7.9 km/h - 230º 3:53 TU"
$ws.Range("N4").Value = "4.2 °C 3:41 TU"
$ws.Range("O4").Value = "7.3 °C"
$ws.Range("E5").Value = "2026-02-26 04:18:46"
$ws.Range("E6").Value = "2026-02-26 04:18:48"
$ws.Range("N6").Value = "8.2 °C 3:43 TU"
$ws.Range("O6").Value = "9.5 °C"
$ws.Range("E7").Value = "2026-02-26 04:18:51"
$ws.Range("N7").Value = "11.1 °C 3:59 TU"
$ws.Range("E8").Value = "2026-02-26 04:18:54"
$ws.Range("Z1").Formula = "=""95%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H8").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("M8").Value = "9.6 °C 3:59 TU"
$ws.Range("N8").Value = "8.8 °C 3:30 TU"
$ws.Range("E9").Value = "2026-02-26 04:18:56"
$ws.Range("N9").Value = "9.9 °C 3:58 TU"
$ws.Range("O9").Value = "10.8 °C"
$ws.Range("E10").Value = "2026-02-26 04:18:59"
$ws.Range("N10").Value = "3.4 °C 3:43 TU"
$ws.Range("O10").Value = "4.3 °C"
$ws.Range("E11").Value = "2026-02-26 04:19:01"
$ws.Range("N11").Value = "1.3 °C 3:59 TU"
$ws.Range("E12").Value = "2026-02-26 04:19:04"
$ws.Range("E13").Value = "2026-02-26 04:19:07"
$ws.Range("J13").Value = "1031.6 hPa"
$ws.Range("L13").Value = "11.2 km/h - 68º 3:41 TU"
$ws.Range("E14").Value = "2026-02-26 04:19:09"
$ws.Range("N14").Value = "8.9 °C 3:44 TU"
$ws.Range("O14").Value = "10.0 °C"
$ws.Range("E15").Value = "2026-02-26 04:19:11"
$ws.Range("N15").Value = "9.2 °C 3:48 TU"
$ws.Range("O15").Value = "10.4 °C"
$ws.Range("E16").Value = "2026-02-26 04:19:14"
$ws.Range("Z1").Formula = "=""48%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H16").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("L16").Value = "20.5 km/h - 278º 3:47 TU"
$ws.Range("E17").Value = "2026-02-26 04:19:16"
$ws.Range("N17").Value = "5.5 °C 3:40 TU"
$ws.Range("E18").Value = "2026-02-26 04:19:19"
$ws.Range("N18").Value = "8.1 °C 3:51 TU"
$ws.Range("E19").Value = "2026-02-26 04:19:22"
$ws.Range("Z1").Formula = "=""68%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H19").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("L19").Value = "5.0 km/h - 300º 3:55 TU"
$ws.Range("M19").Value = "8.0 °C 3:59 TU"
$ws.Range("E20").Value = "2026-02-26 04:19:24"
$ws.Range("Z1").Formula = "=""58%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H20").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("O20").Value = "0.5 °C"
$ws.Range("E21").Value = "2026-02-26 04:19:27"
$ws.Range("J21").Value = "1028.6 hPa"
$ws.Range("N21").Value = "2.9 °C 3:36 TU"
$ws.Range("O21").Value = "4.5 °C"
$ws.Range("E22").Value = "2026-02-26 04:19:29"
$ws.Range("L22").Value = "19.4 km/h - 349º 3:59 TU"
$ws.Range("E23").Value = "2026-02-26 04:19:32"
$ws.Range("Z1").Formula = "=""44%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H23").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("N23").Value = "1.4 °C 3:32 TU"
$ws.Range("O23").Value = "2.4 °C"
$ws.Range("E24").Value = "2026-02-26 04:19:35"
$ws.Range("Z1").Formula = "=""78%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H24").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("J24").Value = "1025.8 hPa"
$ws.Range("N24").Value = "3.1 °C 3:44 TU"
$ws.Range("O24").Value = "7.5 °C"
$ws.Range("E25").Value = "2026-02-26 04:19:37"
$ws.Range("Z1").Formula = "=""40%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H25").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("O25").Value = "3.0 °C"
$ws.Range("E26").Value = "2026-02-26 04:19:39"
$ws.Range("Z1").Formula = "=""45%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H26").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("J26").Value = "1024.8 hPa"
$ws.Range("M26").Value = "8.7 °C 3:44 TU"
$ws.Range("O26").Value = "7.2 °C"
$ws.Range("E27").Value = "2026-02-26 04:19:42"
$ws.Range("N27").Value = "1.5 °C 3:50 TU"
$ws.Range("O27").Value = "2.4 °C"
$ws.Range("E28").Value = "2026-02-26 04:19:44"
$ws.Range("N28").Value = "7.4 °C 3:58 TU"
$ws.Range("O28").Value = "8.6 °C"
$ws.Range("E29").Value = "2026-02-26 04:19:47"
$ws.Range("L29").Value = "9.4 km/h - 341º 3:58 TU"
$ws.Range("N29").Value = "9.1 °C 3:53 TU"
$ws.Range("O29").Value = "9.8 °C"
$ws.Range("E30").Value = "2026-02-26 04:19:50"
$ws.Range("Z1").Formula = "=""100%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H30").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("N30").Value = "10.2 °C 3:57 TU"
$ws.Range("E31").Value = "2026-02-26 04:19:52"
$ws.Range("J31").Value = "1025.6 hPa"
$ws.Range("N31").Value = "9.9 °C 3:58 TU"
$ws.Range("E32").Value = "2026-02-26 04:19:55"
$ws.Range("Z1").Formula = "=""74%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H32").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("N32").Value = "-0.2 °C 3:43 TU"
$ws.Range("O32").Value = "1.5 °C"
$ws.Range("E33").Value = "2026-02-26 04:19:57"
$ws.Range("Z1").Formula = "=""76%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H33").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("J33").Value = "1029.1 hPa"
$ws.Range("N33").Value = "1.3 °C 3:34 TU"
$ws.Range("O33").Value = "2.8 °C"
$ws.Range("E34").Value = "2026-02-26 04:20:00"
$ws.Range("Z1").Formula = "=""51%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H34").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("O34").Value = "2.4 °C"
$ws.Range("E35").Value = "2026-02-26 04:20:03"
$ws.Range("Z1").Formula = "=""30%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H35").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("J35").Value = "1025.2 hPa"
$ws.Range("N35").Value = "7.8 °C 3:49 TU"
$ws.Range("O35").Value = "9.0 °C"
$ws.Range("E36").Value = "2026-02-26 04:20:05"
$ws.Range("J36").Value = "1026.2 hPa"
$ws.Range("E37").Value = "2026-02-26 04:20:07"
$ws.Range("E38").Value = "2026-02-26 04:20:10"
$ws.Range("N38").Value = "6.0 °C 3:56 TU"
$ws.Range("O38").Value = "8.1 °C"
$ws.Range("E39").Value = "2026-02-26 04:20:12"
$ws.Range("Z1").Formula = "=""38%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H39").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("E40").Value = "2026-02-26 04:20:15"
$ws.Range("J40").Value = "1029.4 hPa"
$ws.Range("E41").Value = "2026-02-26 04:20:17"
$ws.Range("O41").Value = "8.1 °C"
$ws.Range("E42").Value = "2026-02-26 04:20:20"
$ws.Range("M42").Value = "9.0 °C 3:47 TU"
$ws.Range("E43").Value = "2026-02-26 04:20:23"
$ws.Range("O43").Value = "3.0 °C"
$ws.Range("E44").Value = "2026-02-26 04:20:25"
$ws.Range("Z1").Formula = "=""59%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H44").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("N44").Value = "-2.2 °C 3:51 TU"
$ws.Range("O44").Value = "-0.1 °C"
$ws.Range("E45").Value = "2026-02-26 04:20:28"
$ws.Range("Z1").Formula = "=""60%"""
$ws.Range("Z1").Copy() | Out-Null
$ws.Range("H45").PasteSpecial($xlPasteValues) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("Z1").ClearContents() | Out-Null
$ws.Range("J45").Value = "1027.2 hPa"
$ws.Range("N45").Value = "4.6 °C 3:40 TU"
$ws.Range("O45").Value = "6.3 °C"
$ws.Range("E46").Value = "2026-02-26 04:20:30"
$ws.Range("N46").Value = "6.4 °C 3:57 TU"
$ws.Range("O46").Value = "7.8 °C"
